# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# Updates to Sheets/Mateus_Profits.xlsx (per-sheet leve profit data)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 460.66666
$ws.Range("I11").Value = 460.66666
$ws.Range("K11").Value = 460.66666
$ws.Range("M11").Value = -320.66666
$ws.Range("H17").Value = 6668381
$ws.Range("J17").Value = 7144622.5
$ws.Range("L17").Value = 21433867.5
$ws.Range("N17").Value = -21434203.5
$ws.Range("H39").Value = 1150.5714
$ws.Range("I39").Value = 263.75
$ws.Range("K39").Value = 791.25
$ws.Range("M39").Value = -495.25
$ws.Range("H106").Value = 1928.4117
$ws.Range("H107").Value = 671.2727
$ws.Range("I107").Value = 749.4737
$ws.Range("K107").Value = 749.4737
$ws.Range("M107").Value = 1170.5263
$ws.Range("H113").Value = 4581.8887
$ws.Range("J113").Value = 5479.4
$ws.Range("L113").Value = 5479.4
$ws.Range("N113").Value = -11987.4
$ws.Range("H137").Value = 1448.037
$ws.Range("I137").Value = 1286.8695
$ws.Range("K137").Value = 3860.6085
$ws.Range("M137").Value = -1310.6085

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1523.3175
$ws.Range("I32").Value = 1550.3771
$ws.Range("K32").Value = 1550.3771
$ws.Range("M32").Value = -1263.3771
$ws.Range("H45").Value = 4266.091
$ws.Range("I45").Value = 2115.875
$ws.Range("J45").Value = 10000
$ws.Range("K45").Value = 2115.875
$ws.Range("L45").Value = 10000
$ws.Range("M45").Value = -1738.875
$ws.Range("N45").Value = -10754
$ws.Range("H46").Value = 7801
$ws.Range("J46").Value = 8263.75
$ws.Range("L46").Value = 8263.75
$ws.Range("N46").Value = -8901.75
$ws.Range("H74").Value = 3013.8262
$ws.Range("I74").Value = 2268.3438
$ws.Range("J74").Value = 4717.7856
$ws.Range("K74").Value = 2268.3438
$ws.Range("L74").Value = 4717.7856
$ws.Range("M74").Value = -1394.3438
$ws.Range("N74").Value = -6465.7856
$ws.Range("H77").Value = 3013.8262
$ws.Range("I77").Value = 2268.3438
$ws.Range("J77").Value = 4717.7856
$ws.Range("K77").Value = 11341.719
$ws.Range("L77").Value = 23588.928
$ws.Range("M77").Value = -6973.719000000001
$ws.Range("N77").Value = -32324.928
$ws.Range("H122").Value = 2950.3635
$ws.Range("I122").Value = 2924.6667
$ws.Range("K122").Value = 8774.000100000001
$ws.Range("M122").Value = -6324.000100000001
$ws.Range("H132").Value = 4398.077
$ws.Range("I132").Value = 3829.394
$ws.Range("K132").Value = 11488.182
$ws.Range("M132").Value = -8958.181999999999
$ws.Range("H135").Value = 108040.664
$ws.Range("J135").Value = 108040.664
$ws.Range("L135").Value = 108040.664
$ws.Range("N135").Value = -118180.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 10738.643
$ws.Range("I58").Value = 8472.5
$ws.Range("J58").Value = 11645.1
$ws.Range("K58").Value = 8472.5
$ws.Range("L58").Value = 11645.1
$ws.Range("M58").Value = -8269.5
$ws.Range("N58").Value = -12051.1
$ws.Range("H74").Value = 42179.2
$ws.Range("J74").Value = 42179.2
$ws.Range("L74").Value = 42179.2
$ws.Range("N74").Value = -43927.2
$ws.Range("H77").Value = 42179.2
$ws.Range("J77").Value = 42179.2
$ws.Range("L77").Value = 126537.6
$ws.Range("N77").Value = -135273.6
$ws.Range("H93").Value = 2375
$ws.Range("I93").Value = 2375
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2375
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -503
$ws.Range("N93").ClearContents()
$ws.Range("H99").Value = 3807.9375
$ws.Range("J99").Value = 3620
$ws.Range("L99").Value = 3620
$ws.Range("N99").Value = -6616
$ws.Range("H126").Value = 3807.9375
$ws.Range("J126").Value = 3620
$ws.Range("L126").Value = 10860
$ws.Range("N126").Value = -15800
$ws.Range("H134").Value = 7742.5557
$ws.Range("I134").Value = 6953.7144
$ws.Range("K134").Value = 20861.1432
$ws.Range("M134").Value = -18326.1432
$ws.Range("H136").Value = 10738.643
$ws.Range("I136").Value = 8472.5
$ws.Range("J136").Value = 11645.1
$ws.Range("K136").Value = 25417.5
$ws.Range("L136").Value = 34935.3
$ws.Range("M136").Value = -22867.5
$ws.Range("N136").Value = -40035.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 645
$ws.Range("I69").Value = 400
$ws.Range("J69").Value = 890
$ws.Range("K69").Value = 1200
$ws.Range("L69").Value = 2670
$ws.Range("M69").Value = -389
$ws.Range("N69").Value = -4292
$ws.Range("H72").Value = 645
$ws.Range("I72").Value = 400
$ws.Range("J72").Value = 890
$ws.Range("K72").Value = 3600
$ws.Range("L72").Value = 8010
$ws.Range("M72").Value = 456
$ws.Range("N72").Value = -16122
$ws.Range("H98").Value = 406.5
$ws.Range("J98").Value = 446.5
$ws.Range("L98").Value = 1339.5
$ws.Range("N98").Value = -4335.5
$ws.Range("H126").Value = 11304.667
$ws.Range("I126").Value = 7609.6665
$ws.Range("K126").Value = 22828.9995
$ws.Range("M126").Value = -17888.9995
$ws.Range("H128").Value = 700328
$ws.Range("I128").Value = 700328
$ws.Range("K128").Value = 2100984
$ws.Range("M128").Value = -2096004
$ws.Range("H140").Value = 435572.66
$ws.Range("I140").Value = 1432.2858
$ws.Range("K140").Value = 4296.857400000001
$ws.Range("M140").Value = 883.1425999999992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 35000
$ws.Range("J59").Value = 35000
$ws.Range("L59").Value = 35000
$ws.Range("N59").Value = -36166
$ws.Range("H93").Value = 36615
$ws.Range("J93").Value = 36615
$ws.Range("L93").Value = 36615
$ws.Range("N93").Value = -40359
$ws.Range("H102").Value = 3986.7
$ws.Range("I102").Value = 3394.6667
$ws.Range("J102").Value = 4874.75
$ws.Range("K102").Value = 3394.6667
$ws.Range("L102").Value = 4874.75
$ws.Range("M102").Value = -1772.6667
$ws.Range("N102").Value = -8118.75
$ws.Range("H122").Value = 4103.1177
$ws.Range("I122").Value = 3608.077
$ws.Range("J122").Value = 5712
$ws.Range("K122").Value = 10824.231
$ws.Range("L122").Value = 17136
$ws.Range("M122").Value = -8374.231
$ws.Range("N122").Value = -22036
$ws.Range("H126").Value = 4826.25
$ws.Range("I126").Value = 4587.143
$ws.Range("K126").Value = 13761.429
$ws.Range("M126").Value = -11291.429
$ws.Range("H132").Value = 2701.2307
$ws.Range("I132").Value = 2676.3333
$ws.Range("K132").Value = 8028.999899999999
$ws.Range("M132").Value = -5498.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3743.7
$ws.Range("I7").Value = 3749.7144
$ws.Range("K7").Value = 3749.7144
$ws.Range("M7").Value = -3637.7144
$ws.Range("H22").Value = 2266.6667
$ws.Range("J22").Value = 400
$ws.Range("L22").Value = 400
$ws.Range("N22").Value = -990
$ws.Range("H27").Value = 2266.6667
$ws.Range("J27").Value = 400
$ws.Range("L27").Value = 400
$ws.Range("N27").Value = -614
$ws.Range("H61").Value = 252251
$ws.Range("I61").Value = 335334.66
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 335334.66
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -335132.66
$ws.Range("N61").Value = -3404
$ws.Range("H113").Value = 252251
$ws.Range("I113").Value = 335334.66
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 335334.66
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -333164.66
$ws.Range("N113").Value = -7340
$ws.Range("H126").Value = 3743.7
$ws.Range("I126").Value = 3749.7144
$ws.Range("K126").Value = 11249.1432
$ws.Range("M126").Value = -8779.143199999999
$ws.Range("H136").Value = 4723.476
$ws.Range("I136").Value = 4271.9443
$ws.Range("K136").Value = 12815.8329
$ws.Range("M136").Value = -10265.8329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4956.8887
$ws.Range("I62").Value = 4732.9165
$ws.Range("K62").Value = 4732.9165
$ws.Range("M62").Value = -4108.9165
$ws.Range("H65").Value = 4956.8887
$ws.Range("I65").Value = 4732.9165
$ws.Range("K65").Value = 23664.5825
$ws.Range("M65").Value = -20544.5825
$ws.Range("H113").Value = 815.5
$ws.Range("I113").Value = 931.8889
$ws.Range("J113").Value = 466.33334
$ws.Range("K113").Value = 2795.6667
$ws.Range("L113").Value = 1399.00002
$ws.Range("M113").Value = -625.6667000000002
$ws.Range("N113").Value = -5739.000019999999
$ws.Range("H126").Value = 6686.364
$ws.Range("I126").Value = 6505.75
$ws.Range("K126").Value = 19517.25
$ws.Range("M126").Value = -17047.25
